$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 6: MethodId=5, Method="ICP-MS", Description="inductively-coupled-plasma mass-spectrometry"
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "ICP-MS"
$ws.Range("C6").Value = "inductively-coupled-plasma mass-spectrometry"

# New row 7: MethodId=6, Method="DIC analyzer" (no Description)
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "DIC analyzer"

# Column C needs to widen to fit the new (longer) content, matching Excel's
# "best fit" column width behavior after the new rows were added.
$ws.Columns.Item(3).ColumnWidth = 108.66666666666667
